# Survey workbook edit:
#  - Column A (RespondentID) on "Survey Responses" (rows 2-20) switches from
#    numeric ids (1..19) to single-letter text codes ("a".."s"), which pulls
#    in 19 new shared strings.
#  - The active sheet / selection moves from "Topic Subjects" (C1:C11) to
#    "Survey Responses" with the cursor parked on A21 (just past the data).

$wb = $excel.ActiveWorkbook

$wsResponses = $wb.Worksheets.Item("Survey Responses")
$wsSubjects  = $wb.Worksheets.Item("Topic Subjects")

# Replace the numeric respondent ids with letter codes a..s.
$letters = @("a","b","c","d","e","f","g","h","i","j","k","l","m","n","o","p","q","r","s")
for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 2
    $wsResponses.Cells.Item($row, 1).Value = $letters[$i]
}

# Move the active sheet / selection to "Survey Responses", cell A21.
$wsResponses.Activate()
$wsResponses.Range("A21").Select()
